# "version final sin errores"
# - Bump the Version value on the Metadata sheet from 0.4.0 to 0.7.0
# - Remove the Jurisdiction / Chile row (row 11) from the Metadata sheet
#   (sheet2 "Include from Ultimo Curso Apr" is left content-wise unchanged;
#   its shared-string indices just shift automatically)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the "Jurisdiction" / "Chile" row entirely, shifting later rows up.
$ws1.Rows.Item(11).Delete()

# Update the Version property value (row 3, column B) from 0.4.0 to 0.7.0.
$ws1.Cells.Item(3, 2).Value = "0.7.0"
